$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source aggregation now only exports a single (most recent, 2014) row
# per selected country. Drop the 2012 and 2013 data rows so the old 2014
# row (currently row 4) shifts up to become row 2.
$ws.Rows(2).Delete()
$ws.Rows(2).Delete()

# Column C's metric was swapped out for a different indicator.
$ws.Range("C1").Value = "0. Crop production index"

# The new crop-production-index value arrived as literal text "98.25" in
# the source feed (not a numeric cell). Force text entry with a leading
# apostrophe, then strip the resulting quote-prefix style back to the
# sheet's plain (unstyled) look so only the value/type changes.
$ws.Range("C2").Value = "'98.25"
$ws.Range("D2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# New trailing column: employment in agriculture (% of total employment).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)"
$ws.Range("I2").Value = 60.16503728435
